$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42624.625648148147
$ws.Range("B4").Value = -60
$ws.Range("C4").Value = 47
$ws.Range("D4").Value = 52
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 100
$ws.Range("G4").Value = 8835
$ws.Range("H4").Value = 5494
$ws.Range("I4").Value = 713
$ws.Range("J4").Value = 91
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = "Noun"
